$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Row 2 (Beta) values ---
$ws.Range("C2").Value = 40.35273080373261
$ws.Range("E2").Value = 0.05768901818751787
$ws.Range("F2").Value = 40.39666407755084
$ws.Range("G2").Value = 40.01469641258531
$ws.Range("H2").Value = 40.77480449357979
$ws.Range("I2").Value = 0.0007661746251621119
$ws.Range("J2").Value = 0.0007164439754845086
$ws.Range("K2").Value = 0.000859645204429976
$ws.Range("L2").Value = 0.05788640087906465
$ws.Range("M2").Value = 0.05752025359672786
$ws.Range("N2").Value = 0.05825604030698193

# --- Update existing Row 3 (Gamma) values ---
$ws.Range("F3").Value = [double]"1.39212999630728e-05"
$ws.Range("G3").Value = [double]"5.46209874913611e-09"
$ws.Range("H3").Value = [double]"3.974521224459411e-05"
$ws.Range("I3").Value = [double]"1.214760637515673e-05"
$ws.Range("J3").Value = [double]"5.098035290843097e-09"
$ws.Range("K3").Value = [double]"3.457176714674844e-05"
$ws.Range("L3").Value = [double]"1.433233372497211e-05"
$ws.Range("M3").Value = [double]"5.64745414735122e-09"
$ws.Range("N3").Value = [double]"4.089762786646532e-05"

# --- Add new Row 4 (Beta + Gamma) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.35273080373261
$ws.Range("D4").Value = 0.0007432820064133916
$ws.Range("E4").Value = 0.05768901818751787
$ws.Range("F4").Value = 40.39667799885081
$ws.Range("G4").Value = 40.01469641804741
$ws.Range("H4").Value = 40.77484423879204
$ws.Range("I4").Value = 0.0007783222315372687
$ws.Range("J4").Value = 0.0007164490735197994
$ws.Range("K4").Value = 0.0008942169715767244
$ws.Range("L4").Value = 0.05790073321278962
$ws.Range("M4").Value = 0.05752025924418201
$ws.Range("N4").Value = 0.0582969379348484

# Copy the number/border/bold formatting from A3 onto the new A4 cell
# (A2/A3 both use the bold-centered style used for the row-index column)
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
